# Auto-generated edit script: refresh cached market-profit figures
# (currentAveragePrice*/LevePrice*/LeveProfit* columns) as produced by the
# scheduled market-data runner. Values mirror a fresh Universalis pull;
# cells that are not applicable for a given leve (e.g. no HQ listings) are
# cleared entirely to match how the source generator omits them.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3413.3333
$ws.Range("I40").Value = 1933.3334
$ws.Range("K40").Value = 1933.3334
$ws.Range("M40").Value = -1758.3334
$ws.Range("H64").Value = 6300
$ws.Range("I64").Value = 6300
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 6300
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -6052
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 6300
$ws.Range("I67").Value = 6300
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 6300
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -5442
$ws.Range("N67").ClearContents()
$ws.Range("H76").Value = 1500
$ws.Range("I76").Value = 1500
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 1500
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -1185
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 1500
$ws.Range("I79").Value = 1500
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 1500
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -408
$ws.Range("N79").ClearContents()
$ws.Range("H113").Value = 6144.4287
$ws.Range("I113").Value = 6667.5
$ws.Range("J113").Value = 3006
$ws.Range("K113").Value = 6667.5
$ws.Range("L113").Value = 3006
$ws.Range("M113").Value = -3413.5
$ws.Range("N113").Value = -9514

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2804.7693
$ws.Range("I45").Value = 3033
$ws.Range("K45").Value = 3033
$ws.Range("M45").Value = -2656
$ws.Range("H61").Value = 2174.25
$ws.Range("I61").Value = 2174.25
$ws.Range("K61").Value = 2174.25
$ws.Range("M61").Value = -1962.25
$ws.Range("H74").Value = 10483.182
$ws.Range("I74").Value = 10531.6
$ws.Range("K74").Value = 10531.6
$ws.Range("M74").Value = -9657.6
$ws.Range("H77").Value = 10483.182
$ws.Range("I77").Value = 10531.6
$ws.Range("K77").Value = 52658
$ws.Range("M77").Value = -48290
$ws.Range("H97").Value = 299
$ws.Range("I97").Value = 299
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 299
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 197
$ws.Range("N97").ClearContents()
$ws.Range("H110").Value = 3892.077
$ws.Range("I110").Value = 2227
$ws.Range("J110").Value = 13050
$ws.Range("K110").Value = 2227
$ws.Range("L110").Value = 13050
$ws.Range("M110").Value = -182
$ws.Range("N110").Value = -17140
$ws.Range("H136").Value = 2174.25
$ws.Range("I136").Value = 2174.25
$ws.Range("K136").Value = 6522.75
$ws.Range("M136").Value = -3972.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 8321.571
$ws.Range("I94").Value = 7056.1113
$ws.Range("J94").Value = 10599.4
$ws.Range("K94").Value = 7056.1113
$ws.Range("L94").Value = 10599.4
$ws.Range("M94").Value = -6605.1113
$ws.Range("N94").Value = -11501.4
$ws.Range("H95").Value = 33414
$ws.Range("J95").Value = 33414
$ws.Range("L95").Value = 33414
$ws.Range("N95").Value = -38906
$ws.Range("H99").Value = 2666.6667
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502
$ws.Range("H101").Value = 57700
$ws.Range("J101").Value = 57700
$ws.Range("L101").Value = 57700
$ws.Range("N101").Value = -64190
$ws.Range("H102").Value = 29999
$ws.Range("I102").Value = 29999
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 29999
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -26754
$ws.Range("N102").ClearContents()
$ws.Range("H103").Value = 71999.5
$ws.Range("J103").Value = 71999.5
$ws.Range("L103").Value = 71999.5
$ws.Range("N103").Value = -74343.5
$ws.Range("H104").Value = 50000
$ws.Range("J104").Value = 50000
$ws.Range("L104").Value = 50000
$ws.Range("N104").Value = -56988
$ws.Range("H105").Value = 6000
$ws.Range("J105").Value = 2000
$ws.Range("L105").Value = 2000
$ws.Range("N105").Value = -5494
$ws.Range("H107").Value = 1225.4
$ws.Range("I107").Value = 1225.4
$ws.Range("K107").Value = 1225.4
$ws.Range("M107").Value = 694.5999999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1542.875
$ws.Range("I16").Value = 1565.6666
$ws.Range("K16").Value = 1565.6666
$ws.Range("M16").Value = -1278.6666
$ws.Range("H31").Value = 1788.8948
$ws.Range("I31").Value = 1790.1
$ws.Range("J31").Value = 1787.5555
$ws.Range("K31").Value = 1790.1
$ws.Range("L31").Value = 1787.5555
$ws.Range("M31").Value = -1495.1
$ws.Range("N31").Value = -2377.5555
$ws.Range("H34").Value = 1788.8948
$ws.Range("I34").Value = 1790.1
$ws.Range("J34").Value = 1787.5555
$ws.Range("K34").Value = 1790.1
$ws.Range("L34").Value = 1787.5555
$ws.Range("M34").Value = -1588.1
$ws.Range("N34").Value = -2191.5555
$ws.Range("H59").Value = 6701.3335
$ws.Range("J59").Value = 10000
$ws.Range("L59").Value = 10000
$ws.Range("N59").Value = -12290
$ws.Range("H60").Value = 16657.334
$ws.Range("I60").Value = 5023.25
$ws.Range("J60").Value = 39925.5
$ws.Range("K60").Value = 5023.25
$ws.Range("L60").Value = 39925.5
$ws.Range("M60").Value = -4512.25
$ws.Range("N60").Value = -40947.5
$ws.Range("H93").Value = 6407
$ws.Range("I93").Value = 6407
$ws.Range("K93").Value = 6407
$ws.Range("M93").Value = -4535
$ws.Range("H105").Value = 4139.1
$ws.Range("I105").Value = 4070.4285
$ws.Range("J105").Value = 4299.3335
$ws.Range("K105").Value = 4070.4285
$ws.Range("L105").Value = 4299.3335
$ws.Range("M105").Value = -2323.4285
$ws.Range("N105").Value = -7793.3335
$ws.Range("H107").Value = 1258.9
$ws.Range("I107").Value = 1424.6666
$ws.Range("K107").Value = 1424.6666
$ws.Range("M107").Value = 495.3334
$ws.Range("H113").Value = 1542.875
$ws.Range("I113").Value = 1565.6666
$ws.Range("K113").Value = 1565.6666
$ws.Range("M113").Value = 604.3334
$ws.Range("H132").Value = 2599.7144
$ws.Range("I132").Value = 1823.4445
$ws.Range("K132").Value = 5470.333500000001
$ws.Range("M132").Value = -2940.333500000001
$ws.Range("H134").Value = 11247.875
$ws.Range("I134").Value = 11659.333
$ws.Range("K134").Value = 34977.999
$ws.Range("M134").Value = -32442.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 845
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 845
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2535
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -4157
$ws.Range("H71").Value = 845
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 845
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 7605
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -15717

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9698.75
$ws.Range("I80").Value = 2931.6667
$ws.Range("K80").Value = 2931.6667
$ws.Range("M80").Value = -1933.6667
$ws.Range("H83").Value = 9698.75
$ws.Range("I83").Value = 2931.6667
$ws.Range("K83").Value = 14658.3335
$ws.Range("M83").Value = -9666.333500000001
$ws.Range("H132").Value = 5249.5
$ws.Range("I132").Value = 5333
$ws.Range("K132").Value = 15999
$ws.Range("M132").Value = -13469

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2755
$ws.Range("I7").Value = 2714.5
$ws.Range("J7").Value = 2998
$ws.Range("K7").Value = 2714.5
$ws.Range("L7").Value = 2998
$ws.Range("M7").Value = -2602.5
$ws.Range("N7").Value = -3222
$ws.Range("H16").Value = 327
$ws.Range("I16").Value = 360
$ws.Range("J16").Value = 195
$ws.Range("K16").Value = 360
$ws.Range("L16").Value = 195
$ws.Range("M16").Value = -190
$ws.Range("N16").Value = -535
$ws.Range("H46").Value = 3843.125
$ws.Range("I46").Value = 1285
$ws.Range("J46").Value = 21750
$ws.Range("K46").Value = 1285
$ws.Range("L46").Value = 21750
$ws.Range("M46").Value = -1097
$ws.Range("N46").Value = -22126
$ws.Range("H68").Value = 2511.8333
$ws.Range("I68").Value = 2558.3635
$ws.Range("K68").Value = 2558.3635
$ws.Range("M68").Value = -1809.3635
$ws.Range("H71").Value = 2511.8333
$ws.Range("I71").Value = 2558.3635
$ws.Range("K71").Value = 12791.8175
$ws.Range("M71").Value = -9047.817499999999
$ws.Range("H122").Value = 7143.625
$ws.Range("I122").Value = 7449.857
$ws.Range("K122").Value = 22349.571
$ws.Range("M122").Value = -19899.571
$ws.Range("H126").Value = 2755
$ws.Range("I126").Value = 2714.5
$ws.Range("J126").Value = 2998
$ws.Range("K126").Value = 8143.5
$ws.Range("L126").Value = 8994
$ws.Range("M126").Value = -5673.5
$ws.Range("N126").Value = -13934

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 32592
$ws.Range("J54").Value = 32592
$ws.Range("L54").Value = 32592
$ws.Range("N54").Value = -33632
$ws.Range("H81").Value = 1854.375
$ws.Range("I81").Value = 1854.375
$ws.Range("K81").Value = 3708.75
$ws.Range("M81").Value = -2647.75
$ws.Range("H84").Value = 1854.375
$ws.Range("I84").Value = 1854.375
$ws.Range("K84").Value = 18543.75
$ws.Range("M84").Value = -13239.75
